$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 updates
$ws.Range("G5").Value = 1.73
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 2.4
$ws.Range("L5").Value = 5.5
$ws.Range("N5").Value = 7.5
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.67
$ws.Range("Z5").Value = 13
$ws.Range("AL5").Value = 41
$ws.Range("AO5").Value = 9.5

# Row 6 update
$ws.Range("N6").Value = 6.8

$wb.Save()
